{"js": "// The edit removes three paragraphs that used to sit right after the\n// \"LOB1255: Hidrologia Aplicada (Requisito fraco)\" paragraph at the end of\n// the document: an empty paragraph, an empty page-break paragraph, and the\n// paragraph holding the \"\u00a9 2020 ... Contact: luizeleno@usp.br ...\" footer\n// text. Locate the anchor paragraph by its text, then delete the three\n// paragraphs that immediately follow it.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst anchorText = \"LOB1255: Hidrologia Aplicada (Requisito fraco)\";\nconst footerText =\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find anchor paragraph: \" + anchorText);\n}\n\n// Sanity-check that the paragraph three slots after the anchor is indeed\n// the footer/copyright paragraph before removing anything (anchor, empty,\n// empty page-break, footer).\nif (paragraphs.items[anchorIndex + 3].text !== footerText) {\n  throw new Error(\"Unexpected document structure near the anchor paragraph.\");\n}\n\n// Removing index (anchorIndex + 1) three times in a row deletes the empty\n// paragraph, the empty page-break paragraph, and the footer paragraph,\n// since each deletion shifts the following paragraphs up by one.\nparagraphs.items[anchorIndex + 1].delete();\nparagraphs.items[anchorIndex + 1].delete();\nparagraphs.items[anchorIndex + 1].delete();\n\nawait context.sync();\n", "ps1": "# The edit removes three paragraphs that used to sit right after the\n# \"LOB1255: Hidrologia Aplicada (Requisito fraco)\" paragraph at the end of\n# the document: an empty paragraph, an empty page-break paragraph, and the\n# paragraph holding the \"\u00a9 2020 ... Contact: luizeleno@usp.br ...\" footer\n# text. Locate the anchor paragraph by its text, then delete the three\n# paragraphs that immediately follow it.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"LOB1255: Hidrologia Aplicada (Requisito fraco)\"\n$footerText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find anchor paragraph: $anchorText\"\n}\n\n# Sanity-check that the paragraph three slots after the anchor is indeed\n# the footer/copyright paragraph before removing anything (anchor, empty,\n# empty page-break, footer).\n$footerIndex = $anchorIndex + 3\n$footerCandidate = $d.Paragraphs.Item($footerIndex).Range.Text.TrimEnd(\"`r\", \"`a\")\nif ($footerCandidate -ne $footerText) {\n    throw \"Unexpected document structure near the anchor paragraph.\"\n}\n\n# Removing the paragraph right after the anchor three times in a row deletes\n# the empty paragraph, the empty page-break paragraph, and the footer\n# paragraph, since each deletion shifts the following paragraphs up by one.\n$d.Paragraphs.Item($anchorIndex + 1).Range.Delete()\n$d.Paragraphs.Item($anchorIndex + 1).Range.Delete()\n$d.Paragraphs.Item($anchorIndex + 1).Range.Delete()\n"}
